$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new Price (D) text (or $null if unchanged), new Volume(1h) (E) text (or $null if unchanged)
$updates = @(
    @(2, '27.083.04', '  +0.43%  '),
    @(3, '1.564.84', '  +0.60%  '),
    @(4, $null, '  +0.62%  '),
    @(5, '210.67', '  +1.70%  '),
    @(6, $null, '  +0.31%  '),
    @(7, $null, '  +0.46%  '),
    @(8, '21.94', '  -0.60%  '),
    @(9, $null, '  +0.09%  '),
    @(10, '0.0598', '  +0.34%  '),
    @(11, '0.0862', '  +0.61%  '),
    @(12, '1.786.38', '  +0.56%  '),
    @(13, '1.556.83', '  +0.06%  '),
    @(14, $null, '  +0.28%  '),
    @(15, $null, '  -0.48%  '),
    @(16, '27.042.88', '  +0.33%  '),
    @(17, '62.02', '  +0.42%  '),
    @(18, '0.0₃0702', '  -0.80%  '),
    @(19, '215.36', '  -1.16%  '),
    @(20, '7.37', '  +0.70%  '),
    @(21, $null, '  +0.51%  '),
    @(22, $null, '  +0.95%  '),
    @(23, '9.21', '  +0.21%  '),
    @(24, '1.94', '  -0.20%  '),
    @(25, '153.90', '  +0.39%  '),
    @(26, '6.60', '  -0.58%  '),
    @(27, '15.04', '  +0.16%  '),
    @(28, $null, '  +1.47%  '),
    @(29, $null, '  +0.44%  '),
    @(30, $null, '  +4.51%  '),
    @(31, $null, '  +0.28%  '),
    @(32, $null, '  +0.16%  '),
    @(33, $null, '  +2.08%  '),
    @(34, '1.440.30', '  +1.51%  '),
    @(35, $null, '  +1.54%  '),
    @(36, $null, '  -0.37%  '),
    @(37, $null, '  +1.88%  '),
    @(38, $null, '  +0.83%  '),
    @(39, '0.531', '  +0.36%  '),
    @(40, $null, '  +2.85%  '),
    @(41, $null, '  -0.16%  '),
    @(42, $null, '  +0.48%  '),
    @(43, '2.35', '  +1.56%  '),
    @(44, $null, '  -0.02%  '),
    @(45, '64.47', '  -0.10%  '),
    @(46, '1.74', '  -0.30%  '),
    @(47, '1.700.09', $null),
    @(48, '86.01', '  -1.50%  '),
    @(49, $null, '  +4.01%  '),
    @(50, $null, '  -0.51%  '),
    @(51, $null, '  -0.22%  ')
)

foreach ($u in $updates) {
    $rowNum = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        # Leading apostrophe forces Excel to store the value as literal text
        # (prevents values like "210.67" or "0.0₃0702" from being auto-converted
        # into numbers / dates). Resetting the style afterwards keeps the cell on
        # its original default formatting.
        $ws.Cells.Item($rowNum, 4).Value = "'" + $dVal
        $ws.Cells.Item($rowNum, 4).Style = "Normal"
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($rowNum, 5).Value = $eVal
    }
}
